# services_charge.xlsx maintenance edit:
#  - remove the "Presentación del servicio" and "Envío" columns
#  - add cell comments documenting each remaining column on row 1
#  - move the active selection to D7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "Presentación del servicio" (E) and "Envío" (H) columns ---
# Deleting E first shifts the old "Envío" column (H) left to G, so it is
# removed second from its new position.
$ws.Range("E1").EntireColumn.Delete()
$ws.Range("G1").EntireColumn.Delete()

# --- Cell comments (Isaac Abensur) describing the bulk-upload format ---
$ws.Range("A1").AddComment("regalalo.pe:`nEl nombre del servicio es obligatorio. Sino se consigna no se registrará en la plataforma.")

$ws.Range("B1").AddComment("Regalalo.pe:`nEl código SKU es único. Si existe entre los productos que haya creado hasta el momento, los datos se actualizarán, en su defecto, se creará un  nuevo registro.")

$ws.Range("C1").AddComment("regalalo.pe:`nEl descuento se consigna del 0 al 99 en la parte entera y como máximo dos digitos en la parte decimal.")

$ws.Range("D1").AddComment("Regalalo.pe:`nEl precio solo acepta números sin símbolos de moneda, sin separación de miles o millares y como máximo dos digitos en la parte decimal.`n")

$ws.Range("E1").AddComment("regalalo.pe:`nLa descripción acepta todo tipo de texto.")

$ws.Range("F1").AddComment("regalalo.pe:`nEl rango de edades a los cuales está dirigido el producto se consigna como sigue:`nedad_mínima,edad_maxima`nDonde:`nedad_mínima: Solo puede ser un valor entero y debe ser menor a la edad máxima`nedad_maxima: Solo puede ser un valor entero y debe ser mayor a la edad mínima")

# --- Restore the active selection recorded in the sheet view ---
$ws.Range("D7").Select()
